# Finish the wireframes and update url map:
#  - fix "-/Home/MyAccount" -> "-/Home/My-Account" on the "URLS continued"
#    slide
#  - add two new sub-paths under it:
#      "-/Home/My-Account/My-Reviews"
#      "-/Home/My-Account/Review-Submission"
#    (each indented with a leading tab, as in the other nested URL rows)

$p = $ppt.ActivePresentation

# Locate the "URLS continued" slide / body placeholder that lists the
# "-/Home/MyAccount" entry instead of hard-coding slide/shape indices.
$targetShape = $null
foreach ($slide in $p.Slides) {
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text.IndexOf("-/Home/MyAccount") -ge 0) {
                $targetShape = $shape
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not find the shape containing '-/Home/MyAccount'"
}

$tr = $targetShape.TextFrame.TextRange

# Rewrite the whole run's text (not just a substring) so PowerPoint
# collapses it back into a single run, keeping the paragraph's existing
# formatting untouched.
$oldUrl = "-/Home/MyAccount"
$newUrl = "-/Home/My-Account"
$fullText = $tr.Text
$startPos = $fullText.IndexOf($oldUrl)
$run = $tr.Characters($startPos + 1, $oldUrl.Length)
$run.Text = $newUrl

# Re-find the paragraph that now reads "-/Home/My-Account" and append two
# sibling paragraphs after it (same paragraph formatting / run formatting
# is inherited automatically from the paragraph we insert after).
$fullText = $tr.Text
$paraCount = $tr.Paragraphs().Count
$accountParaIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $tr.Paragraphs($i, 1)
    # Use StartsWith rather than an exact match: after growing the run's
    # text via Characters(...).Text=, this engine's Paragraphs(i,1).Text
    # can include a trailing paragraph-mark character, so an exact -eq
    # comparison would never hit.
    if ($candidate.Text.StartsWith($newUrl)) {
        $accountParaIndex = $i
    }
}

if ($accountParaIndex -eq -1) {
    throw "Could not re-locate the '-/Home/My-Account' paragraph"
}

$accountPara = $tr.Paragraphs($accountParaIndex, 1)
$tab = [char]9
$newParagraphs = "`r" + $tab + "-/Home/My-Account/My-Reviews" + "`r" + $tab + "-/Home/My-Account/Review-Submission"
$accountPara.InsertAfter($newParagraphs) | Out-Null
